# ModelRuns_Round2.xlsx — copy pathway 4, 1x, 1a and 1b from round 1
# (adds rows 9-15 to the all_runs sheet, with colour-coded fills per pathway
#  and a hyperlink on the Pathway 1x asana-task cell)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New "NoProject" rows (9-11) - same plain styling as existing rows 5-8
# ---------------------------------------------------------------------
$noProjectRows = @(
    @{ Row = 9;  C = "2035_TM160_NGF_r2_NoProject_01" },
    @{ Row = 10; C = "2035_TM160_NGF_r2_NoProject_01_AOCx1.25_v2" },
    @{ Row = 11; C = "2035_TM160_NGF_r2_NoProject_03_pretollcalib" }
)

foreach ($r in $noProjectRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = "NextGenFwys"
    $ws.Range("B$row").Value = 2035
    $ws.Range("B$row").HorizontalAlignment = -4108
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = "NGF_Round2"
    $ws.Range("G$row").Value = "current"
}

# ---------------------------------------------------------------------
# 2. Pathway rows (12-15), copied in from the Round 1 workbook, each
#    colour-coded with a light fill.
# ---------------------------------------------------------------------
$finalBlueprint = '"Final Blueprint runs\Final Blueprint (s24)\BAUS v2.25 - FINAL VERSION"'

$pathwayRows = @(
    @{ Row = 12; C = "2035_TM152_NGF_NP10_Path1a_02"; E = "Pathway 1a";
       F = "Pathway 1a - All Lane Tolling + Transit Double Down";
       H = "NGF_Networks_P1a_AllLaneTolling_ImproveTransit_09";
       I = "Rerun Pathway 1a with new network"; Link = $null; Fill = 13431551 },
    @{ Row = 13; C = "2035_TM152_NGF_NP10_Path1b_02"; E = "Pathway 1b";
       F = "Pathway 1b - All Lane Tolling + Affordable";
       H = "NGF_Networks_P1b_AllLaneTolling_Affordable_04";
       I = "Rerun Pathway 1b with new network"; Link = $null; Fill = 14083579 },
    @{ Row = 14; C = "2035_TM152_NGF_NP10_Path1x_01"; E = "Pathway 1x";
       F = "Pathway 1x - All-lane tolling pricing strategy only";
       H = "NGF_Network_P1x_AllLaneTolling_PricingOnly_01";
       I = "https://app.asana.com/0/1201809392759895/1205309291141002/f";
       Link = "https://app.asana.com/0/1201809392759895/1205309291141002/f"; Fill = 14348258 },
    @{ Row = 15; C = "2035_TM152_NGF_NP10_Path4_02"; E = "Pathway 4";
       F = "Pathway 4 - No New Pricing";
       H = "NGF_Networks_P4_NoNewPricing_03";
       I = "Rerun Pathway 4 with new network"; Link = $null; Fill = 13431551 }
)

# Seed the shared "Hyperlink" font/cell-style from a real hyperlink FIRST
# (row 14's asana link) so it picks up the themed hyperlink colour/size
# rather than a generic default - applying the named style before any
# real hyperlink exists would bake in a plain-RGB font instead.
$linkRow = ($pathwayRows | Where-Object { $_.Link })[0]
$ws.Range("I" + $linkRow.Row).Value = $linkRow.I
$ws.Hyperlinks.Add($ws.Range("I" + $linkRow.Row), $linkRow.Link) | Out-Null

foreach ($r in $pathwayRows) {
    $row = $r.Row

    $ws.Range("A$row").Value = "NextGenFwys"
    $ws.Range("B$row").Value = 2035
    $ws.Range("B$row").HorizontalAlignment = -4108
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = "NGF"
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = "current"
    $ws.Range("H$row").Value = $r.H

    if (-not $r.Link) {
        $ws.Range("I$row").Value = $r.I
        $ws.Range("I$row").Style = "Hyperlink"
    }

    $ws.Range("J$row").Value = $finalBlueprint
    $ws.Range("K$row").Value = "run182"
    $ws.Range("K$row").HorizontalAlignment = -4108

    # colour-code the whole row A:K with the pathway's fill colour
    $ws.Range("A" + $row + ":K" + $row).Interior.Color = $r.Fill
}

# ---------------------------------------------------------------------
# 3. Column widths / autofit for the newly-populated columns
# ---------------------------------------------------------------------
$ws.Columns("E").ColumnWidth = 22.8984375
$ws.Columns("I").ColumnWidth = 59.69921875

# ---------------------------------------------------------------------
# 4. View state - unfreeze the horizontal scroll and move the selection
# ---------------------------------------------------------------------
$ws.Range("D24").Select() | Out-Null

Write-Host "NextGenFwys round-2 rows added"
